$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $escaped = $text.Replace('"', '""')
    $c.Formula = '=TEXT("' + $escaped + '","@")'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

$ws.Range("E2").Value = '2026-02-07 04:47:37'
Set-TextValue "H2" "92%"
$ws.Range("N2").Value = '-2.2 °C 4:29 TU'
$ws.Range("E3").Value = '2026-02-07 04:47:40'
$ws.Range("O3").Value = '-5.6 °C'
$ws.Range("E4").Value = '2026-02-07 04:47:42'
$ws.Range("J4").Value = '1000.9 hPa'
$ws.Range("O4").Value = '11.6 °C'
$ws.Range("E5").Value = '2026-02-07 04:47:45'
$ws.Range("J5").Value = '1001.0 hPa'
$ws.Range("E6").Value = '2026-02-07 04:47:47'
$ws.Range("J6").Value = '1002.7 hPa'
$ws.Range("N6").Value = '11.3 °C 4:02 TU'
$ws.Range("E7").Value = '2026-02-07 04:47:49'
$ws.Range("J7").Value = '1002.4 hPa'
$ws.Range("N7").Value = '7.0 °C 4:01 TU'
$ws.Range("O7").Value = '7.8 °C'
$ws.Range("E8").Value = '2026-02-07 04:47:52'
Set-TextValue "H8" "94%"
$ws.Range("N8").Value = '2.8 °C 4:01 TU'
$ws.Range("O8").Value = '4.3 °C'
$ws.Range("E9").Value = '2026-02-07 04:47:54'
$ws.Range("N9").Value = '0.1 °C 4:28 TU'
$ws.Range("O9").Value = '2.1 °C'
$ws.Range("E10").Value = '2026-02-07 04:47:56'
$ws.Range("E11").Value = '2026-02-07 04:47:59'
$ws.Range("E12").Value = '2026-02-07 04:48:01'
$ws.Range("N12").Value = '8.0 °C 4:29 TU'
$ws.Range("O12").Value = '9.9 °C'
$ws.Range("E13").Value = '2026-02-07 04:48:03'
$ws.Range("E14").Value = '2026-02-07 04:48:05'
Set-TextValue "H14" "78%"
$ws.Range("E15").Value = '2026-02-07 04:48:08'
Set-TextValue "H15" "83%"
$ws.Range("J15").Value = '1001.3 hPa'
$ws.Range("N15").Value = '4.3 °C 4:10 TU'
$ws.Range("O15").Value = '6.9 °C'
$ws.Range("E16").Value = '2026-02-07 04:48:10'
Set-TextValue "H16" "90%"
$ws.Range("N16").Value = '2.0 °C 4:28 TU'
$ws.Range("O16").Value = '3.2 °C'
$ws.Range("E17").Value = '2026-02-07 04:48:12'
Set-TextValue "H17" "98%"
$ws.Range("J17").Value = '1004.5 hPa'
$ws.Range("O17").Value = '3.4 °C'
$ws.Range("E18").Value = '2026-02-07 04:48:15'
$ws.Range("L18").Value = '13.7 km/h - 295º 4:16 TU'
$ws.Range("N18").Value = '-8.5 °C 4:29 TU'
$ws.Range("O18").Value = '-7.0 °C'
$ws.Range("E19").Value = '2026-02-07 04:48:17'
$ws.Range("J19").Value = '1005.7 hPa'
$ws.Range("O19").Value = '4.7 °C'
$ws.Range("E20").Value = '2026-02-07 04:48:19'
Set-TextValue "H20" "85%"
$ws.Range("N20").Value = '-5.6 °C 4:14 TU'
$ws.Range("O20").Value = '-4.5 °C'
$ws.Range("E21").Value = '2026-02-07 04:48:22'
Set-TextValue "H21" "72%"
$ws.Range("J21").Value = '1001.4 hPa'
$ws.Range("N21").Value = '3.5 °C 4:28 TU'
$ws.Range("O21").Value = '7.4 °C'
$ws.Range("E22").Value = '2026-02-07 04:48:24'
Set-TextValue "H22" "94%"
$ws.Range("L22").Value = '14.8 km/h - 3º 4:04 TU'
$ws.Range("M22").Value = '7.4 °C 4:29 TU'
$ws.Range("O22").Value = '5.7 °C'
$ws.Range("E23").Value = '2026-02-07 04:48:26'
$ws.Range("J23").Value = '1001.2 hPa'
$ws.Range("L23").Value = '13.7 km/h - 41º 4:09 TU'
$ws.Range("N23").Value = '6.9 °C 4:24 TU'
$ws.Range("O23").Value = '7.6 °C'
$ws.Range("E24").Value = '2026-02-07 04:48:29'
Set-TextValue "H24" "82%"
$ws.Range("J24").Value = '1000.5 hPa'
$ws.Range("L24").Value = '37.1 km/h - 346º 4:06 TU'
$ws.Range("E25").Value = '2026-02-07 04:48:31'
$ws.Range("J25").Value = '1005.0 hPa'
$ws.Range("O25").Value = '0.6 °C'
$ws.Range("E26").Value = '2026-02-07 04:48:34'
$ws.Range("N26").Value = '-3.1 °C 4:29 TU'
$ws.Range("O26").Value = '-1.4 °C'
$ws.Range("E27").Value = '2026-02-07 04:48:36'
$ws.Range("J27").Value = '1001.0 hPa'
$ws.Range("O27").Value = '8.5 °C'
$ws.Range("E28").Value = '2026-02-07 04:48:38'
Set-TextValue "H28" "87%"
$ws.Range("J28").Value = '1003.6 hPa'
$ws.Range("N28").Value = '1.9 °C 4:21 TU'
$ws.Range("O28").Value = '3.5 °C'
$ws.Range("E29").Value = '2026-02-07 04:48:41'
$ws.Range("N29").Value = '9.2 °C 4:06 TU'
$ws.Range("O29").Value = '11.2 °C'
$ws.Range("E30").Value = '2026-02-07 04:48:43'
$ws.Range("O30").Value = '-4.8 °C'
$ws.Range("E31").Value = '2026-02-07 04:48:46'
$ws.Range("J31").Value = '1005.6 hPa'
$ws.Range("E32").Value = '2026-02-07 04:48:48'
Set-TextValue "H32" "61%"
$ws.Range("J32").Value = '1004.1 hPa'
$ws.Range("K32").Value = '-0.1 MJ/m2'
$ws.Range("O32").Value = '11.5 °C'
$ws.Range("E33").Value = '2026-02-07 04:48:51'
Set-TextValue "H33" "89%"
$ws.Range("N33").Value = '5.5 °C 4:24 TU'
$ws.Range("O33").Value = '7.2 °C'
$ws.Range("E34").Value = '2026-02-07 04:48:53'
$ws.Range("N34").Value = '5.1 °C 4:00 TU'
$ws.Range("O34").Value = '6.6 °C'
$ws.Range("E35").Value = '2026-02-07 04:48:55'
$ws.Range("N35").Value = '-8.2 °C 4:23 TU'
$ws.Range("O35").Value = '-5.1 °C'
$ws.Range("E36").Value = '2026-02-07 04:48:58'
$ws.Range("J36").Value = '1006.3 hPa'
$ws.Range("N36").Value = '4.0 °C 4:02 TU'
$ws.Range("O36").Value = '4.6 °C'

$excel.CutCopyMode = $false
